# Update the QuoteID / date test-data values in the "Input" sheet so the
# QA suite and the staging suite use separate (non-overlapping) order
# references, per commit "Separate suites for qa and stging and message update".
#
# Mapping of old -> new values:
#   Q2  : 51501022 -> 51503454
#   Q3  : 51501023 -> 51503455
#   R3  : 51501024 -> 51503456
#   AD3 : 10-28-2021 -> 11-08-2021
#   Q4  : 51501025 -> 51503457
#
# All of these values must remain stored as *text* (they are ids / dates
# written in a non-ISO format), so each cell's number format is forced to
# "@" (Text) before the new value is assigned; otherwise Excel would
# auto-detect the numeric-looking strings as numbers, or the dash-separated
# date string as an actual date serial.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

$updates = @(
    @{ Cell = "Q2";  Value = "51503454" },
    @{ Cell = "Q3";  Value = "51503455" },
    @{ Cell = "R3";  Value = "51503456" },
    @{ Cell = "AD3"; Value = "11-08-2021" },
    @{ Cell = "Q4";  Value = "51503457" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
